# Daily update at 8 AM UTC
# Appends the next day's row (row 86) to the "Wins Over Time" tracker and
# moves the "latest row" date formatting down onto the new last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous last row (85) used a plain "date only" number format to mark
# it as the latest entry; now that row 86 is the latest entry, row 85 goes
# back to the regular date+time format shared by every other data row.
$ws.Range("A85").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's data.
$ws.Range("A86").Value = 45673
$ws.Range("A86").NumberFormat = "YYYY-MM-DD"
$ws.Range("B86").Value = 203
$ws.Range("C86").Value = 201
$ws.Range("D86").Value = 198
